$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Tom Holland
$ws.Range("A2").Value = "Tom"
$ws.Range("B2").Value = "Sr"
$ws.Range("C2").Value = "Holland"
$ws.Range("E2").Value = "tom123ho"

# Row 3 - Jackie Chan
$ws.Range("A3").Value = "Jackie"
$ws.Range("B3").Value = "Jr"
$ws.Range("C3").Value = "Chan"
$ws.Range("E3").Value = "jack8943ch"

# Row 4 - Star Lord
$ws.Range("A4").Value = "Star"
$ws.Range("B4").Value = "Cool"
$ws.Range("C4").Value = "Lord"
$ws.Range("E4").Value = "star130lord"

# Update selection from C2 to B2
$ws.Range("B2").Select()
